$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.231.13'
$ws.Range("E2").Value = '  -2.68%  '
$ws.Range("D3").Value = '1.701.82'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.07'
$ws.Range("E5").Value = '  -2.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5293'
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2643'
$ws.Range("E8").Value = '  -4.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06566'
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.74'
$ws.Range("E10").Value = '  -4.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07620'
$ws.Range("E11").Value = '  -2.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.561'
$ws.Range("E12").Value = '  -3.11%  '
$ws.Range("D13").Value = '1.731.22'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '1.937.78'
$ws.Range("E14").Value = '  -2.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5705'
$ws.Range("E15").Value = '  -5.20%  '
$ws.Range("E16").Value = '  -3.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.43'
$ws.Range("E17").Value = '  -2.95%  '
$ws.Range("D18").Value = '27.229.44'
$ws.Range("E18").Value = '  -2.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.15'
$ws.Range("E19").Value = '  -4.87%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.658'
$ws.Range("E21").Value = '  -3.91%  '
$ws.Range("E22").Value = '  -5.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.939'
$ws.Range("E23").Value = '  -4.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.35'
$ws.Range("E25").Value = '  -3.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.751'
$ws.Range("E26").Value = '  +5.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1214'
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.233'
$ws.Range("E28").Value = '  -3.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.27'
$ws.Range("E29").Value = '  -4.88%  '
$ws.Range("E30").Value = '  -5.83%  '
$ws.Range("E31").Value = '  -2.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.486'
$ws.Range("E32").Value = '  -5.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.413'
$ws.Range("E33").Value = '  -3.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.637'
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.869'
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.418'
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9440'
$ws.Range("E37").Value = '  -4.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5830'
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01627'
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.853'
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = '1.041.46'
$ws.Range("E42").Value = '  -0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8360'
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.70'
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").Value = '1.844.66'
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("E47").Value = '  -3.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4497'
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.079'
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05238'
$ws.Range("E51").Value = '  -1.56%  '
